# Apply the target edit to the workbook:
# 1. Delete the "Texas Data" worksheet entirely (it is dropped from the
#    workbook, along with its only notes that lived in the shared strings
#    table and its dedicated cell style/font).
# 2. Restore the original (pre-fix) formula in HPEbP!B3, changing it from
#    118/(162+2) back to 118/(162+2+46), which also recalculates every
#    formula that chains off of it across row 3.

# Suppress the "this sheet contains data" confirmation dialog that Excel
# normally raises when deleting a non-empty worksheet.
$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# --- 1. Delete the "Texas Data" sheet ---
$texasSheet = $wb.Worksheets.Item("Texas Data")
$texasSheet.Delete()

# --- 2. Revert the formula in the HPEbP sheet, cell B3 ---
$hpebp = $wb.Worksheets.Item("HPEbP")
$hpebp.Range("B3").Formula = "=118/(162+2+46)"

$excel.DisplayAlerts = $true
